$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.175.58'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '3.600.85'
$ws.Range('E3').Value = '  +1.96%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '603.08'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.12%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '139.32'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.13%  '
$ws.Range('D7').Value = '3.600.37'
$ws.Range('E7').Value = '  +1.95%  '
$ws.Range('E8').Value = '  -0.01%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.501'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.70%  '
$ws.Range('E10').Value = '  +2.49%  '
$ws.Range('E11').Value = '  +4.56%  '
$ws.Range('E12').Value = '  +2.17%  '
$ws.Range('D13').Value = '4.214.51'
$ws.Range('E13').Value = '  +2.08%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '28.39'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +4.70%  '
$ws.Range('E15').Value = '  +2.72%  '
$ws.Range('D16').Value = '3.600.86'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '66.243.60'
$ws.Range('E18').Value = '  +1.75%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '10.16'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.56%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '14.66'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +2.85%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.90'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.08%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '397.50'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.66%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.590'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +3.10%  '
$ws.Range('D24').Value = '3.747.61'
$ws.Range('E24').Value = '  +2.04%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '75.15'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  +6.02%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '8.18'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +4.92%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.64'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +26.66%  '
$ws.Range('E30').Value = '  +6.72%  '
$ws.Range('E31').Value = '  +3.55%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('D33').Value = '3.612.37'
$ws.Range('E33').Value = '  +1.73%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '24.64'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +3.63%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.150'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +4.39%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  +8.75%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.63'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +4.15%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '7.05'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +2.06%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '168.58'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.86%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0844'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +5.36%  '
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('E43').Value = '  +6.66%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '26.30'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('E45').Value = '  +1.46%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '4.56'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +3.33%  '
$ws.Range('E47').Value = '  +0.07%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.72'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +3.35%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '7.01'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +3.35%  '
$ws.Range('D50').Value = '2.459.30'
$ws.Range('E50').Value = '  +3.35%  '
$ws.Range('E51').Value = '  +10.26%  '
